$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Задание" paragraph: expand the task description.
#    The three Find/Replace calls below are applied in an order that keeps
#    each search pattern unambiguous at the moment it runs (replacement (a)
#    inserts the literal text "полупрозрачное ", which is exactly the
#    pattern replacement (b) looks for - so (b) must run first).
# ---------------------------------------------------------------------------

# (b) the run that used to read "полупрозрачное " now introduces the list
#     of characteristics that follows.
$r1 = $d.Content.Find.Execute(
    "полупрозрачное ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "В описании объекта должны присутствовать такие характеристики, как",
    2)

# (a) the adjective phrase that used to start the next sentence is folded
#     into the description of the object right after "представлять собой ".
$r2 = $d.Content.Find.Execute(
    "Объект должен представлять собой ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Объект должен представлять собой полупрозрачное вязкоупругое тело, имеющее форму шара при отсутствии внешних сил. ",
    2)

# (c) the remainder of the old sentence is replaced by the enumeration of
#     characteristics (mass, damping coefficient, transmittance), followed
#     by the original "Объект должен лежать на полу..." text (note the
#     doubled space before "Программа" in the new wording).
$r3 = $d.Content.Find.Execute(
    "вязкоупругое тело, имеющее форму шара при отсутствии внешних сил. Объект должен лежать на полу, имеющем структуру. Программа должна предоставить пользователю интерфейс, дающий возможность изменять объект: растягивать, вдавливать, изменять цвет, степень прозрачности.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " масса тела, коэффициент затухания тела и коэффициент пропускания света, проходящего через данное тело. Объект должен лежать на полу, имеющем структуру.  Программа должна предоставить пользователю интерфейс, дающий возможность изменять объект: растягивать, вдавливать, изменять цвет, степень прозрачности.",
    2)

Write-Output ("task paragraph replacements: b=" + $r1 + " a=" + $r2 + " c=" + $r3)

# ---------------------------------------------------------------------------
# 2) Remove the empty paragraph right after "Оформление курсовой работы:"
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Оформление курсовой работы*") {
        $targetIndex = $i
    }
}
if ($targetIndex -ge 1 -and $targetIndex -lt $d.Paragraphs.Count) {
    $next = $d.Paragraphs($targetIndex + 1)
    $next.Range.Delete()
    Write-Output ("removed empty paragraph after heading at index " + $targetIndex)
}

# ---------------------------------------------------------------------------
# 3) Remove the last (empty) paragraph of the document body.
# ---------------------------------------------------------------------------
if ($d.Paragraphs.Count -gt 0) {
    $last = $d.Paragraphs.Last
    $last.Range.Delete()
    Write-Output "removed trailing empty paragraph"
}
